$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Change the date column (C) number format from the built-in m/d/yyyy
#    style to a custom ISO-ish "yyyy-mm-dd;@" format. Apply to every data
#    row, including the two brand-new ones added below.
# ---------------------------------------------------------------------------
$ws.Range("C2:C49").NumberFormat = "yyyy-mm-dd;@"

# ---------------------------------------------------------------------------
# 2. Fix the "added" dates for the last few existing rows (they were
#    mis-dated 2021-08-08 and should actually be 2021-07-08), and correct
#    row 47's id value.
# ---------------------------------------------------------------------------
$ws.Range("C45").Value = "2021-07-08"
$ws.Range("C46").Value = "2021-07-08"
$ws.Range("A47").Value = 46
$ws.Range("C47").Value = "2021-07-08"

# ---------------------------------------------------------------------------
# 3. Append the two new phishing-sample rows.
# ---------------------------------------------------------------------------
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "msg"
$ws.Range("C48").Value = "2021-07-12"
$ws.Range("D48").Value = "Official"
$ws.Range("E48").Value = "shortened"
$ws.Range("F48").Value = "opportunity"
$ws.Range("G48").Value = "mt"
$ws.Range("H48").Value = "no"
$ws.Range("I48").Value = "get a free covid kit"
$ws.Range("J48").Value = "Government"
$ws.Range("K48").Value = "redirects to https://www.restaurant-apron.at/wp-admin/network/-/"

$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "msg"
$ws.Range("C49").Value = "2021-07-12"
$ws.Range("D49").Value = "SocialMedia"
$ws.Range("E49").Value = "shortened"
$ws.Range("F49").Value = "delivery"
$ws.Range("G49").Value = "mt"
$ws.Range("H49").Value = "no"
$ws.Range("I49").Value = "delivery payment"
$ws.Range("J49").Value = "MaltaPost"

# ---------------------------------------------------------------------------
# 4. Re-point the frozen header pane / selection so the view now shows the
#    newly added rows, mirroring the scrolled state captured in the author's
#    workbook.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("A49").Select()
